$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The paragraph right after "Planejamento:" used to describe the software
#    stack (HTML5/CSS3/JS/MySQL/Azure). Its text is replaced by a new
#    paragraph describing the login/registration system requirement.
# ---------------------------------------------------------------------------
$oldPlanejamentoText = "O software será desenvolvido em HTML5, CSS3 e Javascript. Será utilizada uma API em Javascript para gerenciar os dados de cadastro e login, os dados serão salvos em um banco de dados desenvolvido em MySQL/SQL server na nuvem Azure."
$newPlanejamentoText = "O site deverá conter um sistema de cadastro e login conectado com o banco de dados SQL Server se hospedado na nuvem e MySQL se hospedado localmente, com isso ele deverá ter pelo menos uma página para essa funcionalidade. Fora a pagina home deverá também ser criada uma área exclusiva para usuários logados no sistema."

$planejamentoPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $oldPlanejamentoText) {
        $planejamentoPara = $p
        break
    }
}

# Collapse the paragraph (which is made up of several runs) down to a single
# run holding the new text, keeping the formatting of the first run.
$r = $planejamentoPara.Range
$r.End = $r.End - 1   # exclude the paragraph mark
$r.Text = $newPlanejamentoText

# ---------------------------------------------------------------------------
# 2) Insert a new paragraph right after the "Escopo:" heading, moving the old
#    software-stack description text there (it now documents the tech stack
#    actually used, instead of the planning description).
# ---------------------------------------------------------------------------
$qualquerPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Qualquer sistema com acesso")) {
        $qualquerPara = $p
        break
    }
}

$qStart = $qualquerPara.Range.Start
$insertRange = $d.Range($qStart, $qStart)
$insertRange.InsertParagraphBefore()
$newRange = $d.Range($qStart, $qStart)
$newRange.Text = $oldPlanejamentoText

# ---------------------------------------------------------------------------
# 3) Append a clarifying sentence about nodejs to the "Restrições" paragraph.
# ---------------------------------------------------------------------------
$restricoesContentPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "O site deverá ser hospedado localmente.") {
        $restricoesContentPara = $p
        break
    }
}

$contentRange = $restricoesContentPara.Range
$contentRange.End = $contentRange.End - 1
$contentRange.InsertAfter(" Para o usuário ter acesso ao sistema de cadastro e login deverá usar o nodejs.")
